# Refresh cryptos list snapshot: update Price (column D) and
# Volume(1h) (column E) text values for each coin row, matching the
# latest scrape. Row 37/38 also swap which coin (MXToken vs
# RenderToken) occupies that rank position, so Coin (B), Link (C),
# Price (D) and Volume(1h) (E) are all rewritten for those two rows.
#
# All of these values are plain text in the original workbook (no
# leading apostrophe is stored in the saved value), so we prefix each
# assignment with a literal apostrophe to force Excel to keep them as
# text instead of auto-converting number-like strings (e.g. "1.005",
# "14.60", "6.303") into numeric values and losing formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.621.99'
$ws.Range('E2').Value = '''  -2.38%  '
$ws.Range('D3').Value = '''1.860.87'
$ws.Range('E3').Value = '''  -2.53%  '
$ws.Range('E4').Value = '''  +0.13%  '
$ws.Range('D5').Value = '''291.27'
$ws.Range('E5').Value = '''  -5.33%  '
$ws.Range('D6').Value = '''1.005'
$ws.Range('E6').Value = '''  +0.26%  '
$ws.Range('D7').Value = '''0.5249'
$ws.Range('E7').Value = '''  -1.57%  '
$ws.Range('D8').Value = '''0.3697'
$ws.Range('E8').Value = '''  -3.22%  '
$ws.Range('D9').Value = '''0.07100'
$ws.Range('E9').Value = '''  -2.58%  '
$ws.Range('D10').Value = '''21.11'
$ws.Range('E10').Value = '''  -4.35%  '
$ws.Range('E11').Value = '''  -2.74%  '
$ws.Range('D12').Value = '''0.08080'
$ws.Range('E12').Value = '''  -1.49%  '
$ws.Range('D13').Value = '''1.924.42'
$ws.Range('E13').Value = '''  +63.48%  '
$ws.Range('D14').Value = '''91.26'
$ws.Range('E14').Value = '''  -4.73%  '
$ws.Range('D15').Value = '''5.233'
$ws.Range('E15').Value = '''  -2.11%  '
$ws.Range('D16').Value = '''1.003'
$ws.Range('E16').Value = '''  -0.11%  '
$ws.Range('D17').Value = '''14.60'
$ws.Range('E17').Value = '''  -1.46%  '
$ws.Range('D18').Value = '''0.000008416'
$ws.Range('E18').Value = '''  -2.77%  '
$ws.Range('D19').Value = '''1.003'
$ws.Range('E19').Value = '''  +0.12%  '
$ws.Range('D20').Value = '''26.675.56'
$ws.Range('E20').Value = '''  -2.32%  '
$ws.Range('D21').Value = '''4.922'
$ws.Range('E21').Value = '''  -2.52%  '
$ws.Range('D22').Value = '''10.54'
$ws.Range('E22').Value = '''  -2.66%  '
$ws.Range('D23').Value = '''6.303'
$ws.Range('D24').Value = '''144.71'
$ws.Range('E24').Value = '''  -3.54%  '
$ws.Range('D25').Value = '''2.227'
$ws.Range('E25').Value = '''  -2.72%  '
$ws.Range('D26').Value = '''1.736'
$ws.Range('E26').Value = '''  -0.61%  '
$ws.Range('D27').Value = '''17.82'
$ws.Range('E27').Value = '''  -2.45%  '
$ws.Range('D28').Value = '''112.89'
$ws.Range('E28').Value = '''  -3.51%  '
$ws.Range('D29').Value = '''4.642'
$ws.Range('E29').Value = '''  -3.84%  '
$ws.Range('D30').Value = '''4.555'
$ws.Range('E30').Value = '''  -5.39%  '
$ws.Range('D31').Value = '''0.09019'
$ws.Range('E31').Value = '''  -2.96%  '
$ws.Range('D32').Value = '''0.7889'
$ws.Range('E32').Value = '''  -5.94%  '
$ws.Range('D33').Value = '''0.04945'
$ws.Range('E33').Value = '''  -2.33%  '
$ws.Range('E34').Value = '''  -1.92%  '
$ws.Range('D35').Value = '''1.150'
$ws.Range('E35').Value = '''  -6.14%  '
$ws.Range('D36').Value = '''0.5896'
$ws.Range('E36').Value = '''  +2.41%  '
$ws.Range('B37').Value = '''RenderToken'
$ws.Range('C37').Value = '''https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = '''2.593'
$ws.Range('E37').Value = '''  -3.85%  '
$ws.Range('B38').Value = '''MXToken'
$ws.Range('C38').Value = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '''3.173'
$ws.Range('E38').Value = '''  -5.54%  '
$ws.Range('D39').Value = '''0.01919'
$ws.Range('E39').Value = '''  -4.40%  '
$ws.Range('D40').Value = '''1.048'
$ws.Range('E40').Value = '''  -2.70%  '
$ws.Range('D41').Value = '''0.5144'
$ws.Range('E41').Value = '''  +4.38%  '
$ws.Range('D42').Value = '''6.411'
$ws.Range('E42').Value = '''  -2.26%  '
$ws.Range('D43').Value = '''114.25'
$ws.Range('E43').Value = '''  -2.66%  '
$ws.Range('D44').Value = '''8.534'
$ws.Range('E44').Value = '''  -8.17%  '
$ws.Range('D45').Value = '''0.1468'
$ws.Range('E45').Value = '''  -3.71%  '
$ws.Range('D46').Value = '''1.004'
$ws.Range('E46').Value = '''  +0.22%  '
$ws.Range('D47').Value = '''9.878'
$ws.Range('E47').Value = '''  -2.73%  '
$ws.Range('D48').Value = '''1.605'
$ws.Range('E48').Value = '''  -2.00%  '
$ws.Range('D49').Value = '''36.87'
$ws.Range('E49').Value = '''  -4.46%  '
$ws.Range('D50').Value = '''0.06017'
$ws.Range('E50').Value = '''  -2.01%  '
$ws.Range('D51').Value = '''61.81'
$ws.Range('E51').Value = '''  -2.72%  '
